$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the 钱 (E) and 粮 (F) weights for rows 20-23 (中古/启蒙/工业/现代 eras)
$ws.Range("E20").Value = 1.5
$ws.Range("F20").Value = 1

$ws.Range("E21").Value = 1.55
$ws.Range("F21").Value = 1

$ws.Range("E22").Value = 1.6
$ws.Range("F22").Value = 1

$ws.Range("E23").Value = 1.6
$ws.Range("F23").Value = 0.9

# Update the view to reflect the scrolled/selected state from the edit session
$ws.Activate()
$ws.Range("E22").Select()

$wb.Save()
